$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "66.193.90"
$ws.Range("E2").Value = "  +7.40%  "

$ws.Range("D3").Value = "3.016.64"
$ws.Range("E3").Value = "  +4.34%  "

$ws.Range("E4").Value = "  -0.08%  "

Set-TextValue "D5" "585.69"
$ws.Range("E5").Value = "  +3.49%  "

Set-TextValue "D6" "156.27"
$ws.Range("E6").Value = "  +9.64%  "

Set-TextValue "D7" "0.999"

$ws.Range("D8").Value = "3.010.92"
$ws.Range("E8").Value = "  +4.23%  "

Set-TextValue "D9" "0.518"
$ws.Range("E9").Value = "  +3.49%  "

Set-TextValue "D10" "6.96"
$ws.Range("E10").Value = "  +0.99%  "

Set-TextValue "D11" "0.156"
$ws.Range("E11").Value = "  +7.58%  "

Set-TextValue "D12" "0.453"
$ws.Range("E12").Value = "  +6.22%  "

Set-TextValue "D13" "0.0000252"
$ws.Range("E13").Value = "  +9.64%  "

Set-TextValue "D14" "34.69"
$ws.Range("E14").Value = "  +9.80%  "

$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("D16").Value = "66.120.05"
$ws.Range("E16").Value = "  +7.30%  "

$ws.Range("D17").Value = "3.514.87"
$ws.Range("E17").Value = "  +4.25%  "

$ws.Range("E18").Value = "  +7.29%  "

$ws.Range("D19").Value = "3.011.28"
$ws.Range("E19").Value = "  +3.82%  "

Set-TextValue "D20" "464.19"
$ws.Range("E20").Value = "  +8.18%  "

Set-TextValue "D21" "13.85"
$ws.Range("E21").Value = "  +6.98%  "

$ws.Range("E22").Value = "  +5.41%  "

$ws.Range("E23").Value = "  +8.89%  "

Set-TextValue "D24" "81.99"
$ws.Range("E24").Value = "  +3.99%  "

Set-TextValue "D25" "2.27"
$ws.Range("E25").Value = "  +13.56%  "

Set-TextValue "D26" "12.53"
$ws.Range("E26").Value = "  +5.84%  "

Set-TextValue "D27" "10.69"
$ws.Range("E27").Value = "  +8.15%  "

$ws.Range("E28").Value = "  -0.08%  "

Set-TextValue "D29" "8.02"
$ws.Range("E29").Value = "  +14.72%  "

Set-TextValue "D30" "2.37"
$ws.Range("E30").Value = "  +16.78%  "

$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("E32").Value = "  +5.44%  "

$ws.Range("E33").Value = "  +5.84%  "

Set-TextValue "D34" "27.05"
$ws.Range("E34").Value = "  +6.52%  "

$ws.Range("E35").Value = "  -0.18%  "

Set-TextValue "D36" "0.995"
$ws.Range("E36").Value = "  +3.95%  "

$ws.Range("E37").Value = "  +9.11%  "

Set-TextValue "D38" "2.18"
$ws.Range("E38").Value = "  +13.60%  "

Set-TextValue "D39" "3.07"
$ws.Range("E39").Value = "  +10.40%  "

Set-TextValue "D40" "49.34"
$ws.Range("E40").Value = "  +1.18%  "

$ws.Range("E41").Value = "  +9.03%  "

$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D42" "0.304"
$ws.Range("E42").Value = "  +14.53%  "

$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D43" "43.98"
$ws.Range("E43").Value = "  +12.10%  "

Set-TextValue "D44" "8.46"
$ws.Range("E44").Value = "  +4.16%  "

Set-TextValue "D45" "394.14"
$ws.Range("E45").Value = "  +15.09%  "

$ws.Range("D46").Value = "2.797.64"
$ws.Range("E46").Value = "  +4.47%  "

$ws.Range("E47").Value = "  +6.04%  "

Set-TextValue "D48" "134.24"
$ws.Range("E48").Value = "  +1.45%  "

$ws.Range("E49").Value = "  -0.04%  "

Set-TextValue "D50" "23.61"
$ws.Range("E50").Value = "  +10.27%  "

$ws.Range("E51").Value = "  +4.07%  "

